$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> column -> new value, derived from the commit diff.
$changes = @{
    2 = @{ D='57.058.42'; E='  +0.75%  ' }
    3 = @{ D='2.445.51'; E='  -1.82%  ' }
    4 = @{ E='  -0.09%  ' }
    5 = @{ D='487.98'; E='  -0.83%  ' }
    6 = @{ D='155.41'; E='  +2.84%  ' }
    7 = @{ D='0.619'; E='  +20.14%  ' }
    8 = @{ D='0.997'; E='  +0.01%  ' }
    9 = @{ D='2.440.46'; E='  -2.40%  ' }
    10 = @{ D='0.0995'; E='  +1.29%  ' }
    11 = @{ D='5.65'; E='  -2.02%  ' }
    12 = @{ D='0.336'; E='  -0.12%  ' }
    13 = @{ E='  +0.77%  ' }
    14 = @{ D='2.860.54'; E='  -2.24%  ' }
    15 = @{ D='57.167.85'; E='  +0.86%  ' }
    16 = @{ D='20.78'; E='  -2.06%  ' }
    17 = @{ E='  -2.20%  ' }
    18 = @{ D='2.437.81'; E='  -2.36%  ' }
    19 = @{ D='4.75'; E='  +4.87%  ' }
    20 = @{ D='327.82'; E='  +2.00%  ' }
    21 = @{ D='9.97'; E='  -2.99%  ' }
    22 = @{ D='1.00'; E='  +0.20%  ' }
    23 = @{ D='5.90'; E='  -0.05%  ' }
    24 = @{ E='  -1.13%  ' }
    25 = @{ E='  -0.33%  ' }
    26 = @{ E='  -0.34%  ' }
    27 = @{ D='0.158'; E='  -2.90%  ' }
    28 = @{ D='2.536.51'; E='  -2.37%  ' }
    29 = @{ D='7.25'; E='  -5.49%  ' }
    30 = @{ D='0.0₃0784'; E='  -2.16%  ' }
    31 = @{ E='  +0.02%  ' }
    32 = @{ D='18.75'; E='  +2.06%  ' }
    33 = @{ D='148.83'; E='  -1.03%  ' }
    34 = @{ E='  +0.26%  ' }
    35 = @{ D='5.31'; E='  +2.04%  ' }
    36 = @{ E='  -1.17%  ' }
    37 = @{ E='  -2.06%  ' }
    38 = @{ D='0.852'; E='  -2.63%  ' }
    39 = @{ E='  +10.66%  ' }
    40 = @{ D='34.23'; E='  +0.87%  ' }
    41 = @{ D='1.37'; E='  -1.27%  ' }
    42 = @{ D='3.52'; E='  +0.29%  ' }
    43 = @{ D='0.997'; E='  +0.23%  ' }
    44 = @{ D='0.596'; E='  -2.70%  ' }
    45 = @{ D='268.70'; E='  +1.79%  ' }
    46 = @{ D='0.0533'; E='  -4.58%  ' }
    47 = @{ B='WhiteBITCoin'; C='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; D='10.21'; E='  -0.14%  ' }
    48 = @{ B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='4.69'; E='  -3.76%  ' }
    49 = @{ E='  -0.91%  ' }
    50 = @{ D='17.53'; E='  -1.55%  ' }
    51 = @{ D='1.857.24'; E='  -2.37%  ' }
}

foreach ($row in $changes.Keys) {
    $cells = $changes[$row]
    foreach ($col in $cells.Keys) {
        $value = $cells[$col]
        $cell = $ws.Range("$col$row")
        if ($col -eq "D") {
            # Column D holds numeric-looking text (prices). Prefix with an
            # apostrophe so Excel stores it as text instead of coercing it
            # into a number, then reset the style so no stray quote-prefix
            # formatting is left attached to the cell.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
